# Update cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new values look numeric,
# so Excel stores them as text strings (matching original inlineStr formatting)
# rather than re-interpreting them as numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"

# Set D-column (Price) values
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "230.93"
$ws.Range("D7").Value = "58.07"
$ws.Range("D9").Value = "0.387"
$ws.Range("D10").Value = "0.0808"
$ws.Range("D13").Value = "14.60"
$ws.Range("D14").Value = "20.72"
$ws.Range("D15").Value = "0.752"
$ws.Range("D16").Value = "5.27"
$ws.Range("D19").Value = "6.26"
$ws.Range("D20").Value = "70.00"
$ws.Range("D22").Value = "225.44"
$ws.Range("D23").Value = "1.00"
$ws.Range("D26").Value = "9.33"
$ws.Range("D27").Value = "165.97"
$ws.Range("D29").Value = "19.15"
$ws.Range("D34").Value = "4.59"
$ws.Range("D37").Value = "5.99"
$ws.Range("D39").Value = "1.00"
$ws.Range("D40").Value = "98.78"
$ws.Range("D41").Value = "0.0219"
$ws.Range("D43").Value = "0.0953"
$ws.Range("D44").Value = "16.82"
$ws.Range("D2").Value = "38.151.50"
$ws.Range("D3").Value = "2.068.96"
$ws.Range("D12").Value = "2.374.69"
$ws.Range("D17").Value = "2.068.96"
$ws.Range("D18").Value = "38.087.05"
$ws.Range("D42").Value = "1.484.28"
$ws.Range("D51").Value = "2.255.47"

# Set E-column (Volume 1h) values
$ws.Range("E2").Value = "  +2.69%  "
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("E7").Value = "  +5.22%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +1.63%  "
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("E18").Value = "  +2.73%  "
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("E25").Value = "  +2.98%  "
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("E28").Value = "  +5.95%  "
$ws.Range("E29").Value = "  +1.90%  "
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  +1.38%  "
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  +3.26%  "
$ws.Range("E35").Value = "  +7.85%  "
$ws.Range("E37").Value = "  +10.53%  "
$ws.Range("E38").Value = "  +4.53%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  +3.43%  "
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +2.58%  "
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("E45").Value = "  +3.29%  "
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("E47").Value = "  +15.56%  "
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("E51").Value = "  +2.14%  "
